# Added param values for sl_bt_scanner_set_mode() and sl_bt_scanner_set_timing().
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F9 needs to end up sharing the same style as F7/F10 (the "note" style), so copy
# that formatting across before filling in the text.
$ws.Range("F10").Copy()
$ws.Range("F9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("F7").Value = "Set phy to 1M and mode to passive scanning"
$ws.Range("F8").Value = "Set phy to 1M, scan interval to 50ms and scan window to 25ms"
$ws.Range("F9").Value = "See the assignment document for these values"

# F21: the middle run "handle_ble_event()" was specially colored - remove that
# highlight so the whole note renders in the normal text color.
$text = "Do not build your state machine into your  handle_ble_event() function! Instead build a new state machine for A7. "
$idx = $text.IndexOf("handle_ble_event()")
$chars = $ws.Range("F21").Characters($idx + 1, 18)
$chars.Font.Color = 0
